# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect the latest generated data (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 727
$ws1.Range("F3").Value = 39
$ws1.Range("F4").Value = 244
$ws1.Range("F5").Value = 2714
$ws1.Range("F7").Value = 3706
$ws1.Range("F8").Value = 470
$ws1.Range("F9").Value = 935
$ws1.Range("F10").Value = 13

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 727
$ws4.Range("F3").Value = 39
$ws4.Range("F5").Value = 244
$ws4.Range("F6").Value = 2714
$ws4.Range("F8").Value = 3706
$ws4.Range("F9").Value = 470
$ws4.Range("F10").Value = 935
$ws4.Range("F11").Value = 13
